$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 153292
$ws.Range("C4").Value = 144840
$ws.Range("C7").Value = 5.51
$ws.Range("C8").Value = 63.86
